# Pair Programming Log - Sprint 5
# 1) Move the stray "_GoBack" bookmark off the title paragraph.
# 2) Fix the recorded Activity text.
# 3) Re-create the "_GoBack" bookmark in its own paragraph at the end of
#    the document (after the existing trailing empty paragraph).

$d = $word.ActiveDocument

# --- 1) Strip the bookmarkStart/bookmarkEnd pair out of the title paragraph ---
# The title paragraph ("Pair Programming Log") currently also carries the
# _GoBack bookmark. Clearing the paragraph's range removes its run AND any
# bookmarks anchored inside it, then we simply retype the title text.
$titlePara = $d.Paragraphs(1)
$titlePara.Range.Delete()
$d.Paragraphs(1).Range.Text = "Pair Programming Log"

# --- 2) Correct the Activity cell text ---
$d.Content.Find.Execute("Register to server, Sell shares", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Automated Capability", 2)

# --- 3) Append a new paragraph holding the _GoBack bookmark at the end ---
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$bookmarkPara = $d.Paragraphs($d.Paragraphs.Count)
$bookmarkPara.Range.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/></w:p>")
